# edit.ps1 -- apply the commit's changes via PowerPoint COM-interop
#
# 1) Slide 6's table switches from the deck's custom "Table_0" style to the
#    built-in PowerPoint table style {BBC9AC50-8D29-4E55-BC21-2625E0227631}.
# 2) The presentation's theme colour scheme (theme1.xml, the slide master's
#    theme) is swapped from the "Integral" palette to the standard Office
#    theme palette.

$p = $ppt.ActivePresentation

# --- helper: pack R,G,B (0-255) into the BGR-packed long that PowerPoint's
#     RGB color properties expect -----------------------------------------
function Pack-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Table style on slide 6 ---------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BBC9AC50-8D29-4E55-BC21-2625E0227631}")
    }
}

# --- 2) Theme colour scheme: Integral -> Office Theme ----------------------
# Office theme standard palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    @(0x00, 0x00, 0x00),   # dk1      000000
    @(0xFF, 0xFF, 0xFF),   # lt1      FFFFFF
    @(0x44, 0x54, 0x6A),   # dk2      44546A
    @(0xE7, 0xE6, 0xE6),   # lt2      E7E6E6
    @(0x5B, 0x9B, 0xD5),   # accent1  5B9BD5
    @(0xED, 0x7D, 0x31),   # accent2  ED7D31
    @(0xA5, 0xA5, 0xA5),   # accent3  A5A5A5
    @(0xFF, 0xC0, 0x00),   # accent4  FFC000
    @(0x44, 0x72, 0xC4),   # accent5  4472C4
    @(0x70, 0xAD, 0x47),   # accent6  70AD47
    @(0x05, 0x63, 0xC1),   # hlink    0563C1
    @(0x95, 0x4F, 0x72)    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $triplet = $officeColors[$i - 1]
    $themeColors.Colors($i).RGB = Pack-RGB $triplet[0] $triplet[1] $triplet[2]
}
